$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 holds account 004204344 / CLINEO / 60000 — delete it entirely,
# shifting the rows below up (matches the diff removing that <x:row>).
$ws.Rows.Item(4).Delete()
